# "Section Hero and Projects completed"
# Updates the works_at mini-table on the "Members" sheet (rows 19-22):
#  - adds a "title" column (E) and a generated SQL "INSERT INTO works_at
#    VALUES" column (I) driven by a formula
#  - reorders / updates the three works_at rows so the newest role
#    (Cathay Bank Intern, eid 3) leads, followed by the Metropolitan Bank
#    Branch Operations Officer (eid 1, now ended 2/2025) and the earlier
#    Relationship Banker stint (eid 2)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Members")

# ---- header row (19): add "title" column ----
$ws.Range("E19").Value2 = "title"

# ---- row 20: Cathay Bank (eid 3) - Intern, 5/2024 to Current ----
$ws.Range("B20").Value2 = 3
$ws.Range("C20").ClearFormats()
$ws.Range("C20").Value2 = "'5/2024"
$ws.Range("C20").NumberFormat = "mmm-yy"
$ws.Range("D20").ClearFormats()
$ws.Range("D20").Value2 = "Current"
$ws.Range("E20").Value2 = "Intern"

# ---- row 21: Metropolitan Bank (eid 1) - Branch Operations Officer, 10/2019 to 2/2025 ----
$ws.Range("B21").Value2 = 1
$ws.Range("C21").ClearFormats()
$ws.Range("C21").Value2 = "'10/2019"
$ws.Range("D21").ClearFormats()
$ws.Range("D21").Value2 = "'2/2025"
$ws.Range("E21").Value2 = "Branch Operations Officer"

# ---- row 22: Metropolitan Bank (eid 2) - Relationship Banker, 1/2018 to 9/2019 ----
$ws.Range("B22").Value2 = 2
$ws.Range("C22").ClearFormats()
$ws.Range("C22").Value2 = "'1/2018"
$ws.Range("D22").ClearFormats()
$ws.Range("D22").Value2 = "'9/2019"
$ws.Range("E22").Value2 = "Relationship Banker"

# ---- header row (19): label for the generated SQL column ----
$ws.Range("H19").Value2 = "INSERT INTO works_at VALUES"

# ---- column I: generated "INSERT INTO works_at VALUES" row strings ----
$ws.Range("I20").Formula = "=""(""&A20&"", ""&B20&"", '""&C20&""', '""&D20&""', '""&E20&""'),"""
$ws.Range("I21:I22").Formula = "=""(""&A21&"", ""&B21&"", '""&C21&""', '""&D21&""', '""&E21&""'),"""

# ---- view: scroll so the updated block is visible, select the new range ----
$ws.Range("E19:I22").Select()

# ---- page setup: printed in portrait orientation ----
$ws.PageSetup.Orientation = 1
